$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"6.236796896606509e-11"
$ws.Range("C2").Value = [double]"1.352317117787379e-10"
$ws.Range("D2").Value = [double]"2.52885433200074e-09"
$ws.Range("E2").Value = [double]"-8.871129415879855e-09"
$ws.Range("B3").Value = [double]"1.352322444830524e-10"
$ws.Range("C3").Value = [double]"6.236833623425707e-11"
$ws.Range("D3").Value = [double]"-8.871128625095189e-09"
$ws.Range("E3").Value = [double]"2.528855076463291e-09"
$ws.Range("B5").Value = [double]"-8.694295182898798e-10"
$ws.Range("C5").Value = [double]"8.69400082546891e-10"
$ws.Range("D5").Value = [double]"-3.313802661667689e-08"
$ws.Range("E5").Value = [double]"3.313798571457936e-08"
$ws.Range("B6").Value = [double]"-1.876790699643002e-08"
$ws.Range("C6").Value = [double]"-1.876788342162414e-08"
$ws.Range("D6").Value = [double]"-3.652095439094759e-06"
$ws.Range("E6").Value = [double]"-3.652095451348261e-06"
$ws.Range("B7").Value = [double]"-3.657621202889199e-08"
$ws.Range("C7").Value = [double]"3.654057730101874e-08"
$ws.Range("D7").Value = [double]"-6.708925808215793e-06"
$ws.Range("E7").Value = [double]"6.708890173479965e-06"
$ws.Range("B8").Value = [double]"3.99963251297649e-08"
$ws.Range("C8").Value = [double]"8.97363737133241e-13"
$ws.Range("D8").Value = [double]"-2.12294181165438e-10"
$ws.Range("E8").Value = [double]"1.751416719416877e-10"
$ws.Range("B9").Value = [double]"-1.963187452441503e-12"
$ws.Range("C9").Value = [double]"4.752061629249763e-13"
$ws.Range("D9").Value = [double]"-1.543888872019434e-10"
$ws.Range("E9").Value = [double]"1.019253803672413e-10"
$ws.Range("B10").Value = [double]"-1.131215241093897e-13"
$ws.Range("C10").Value = [double]"6.565920595550118e-15"
$ws.Range("D10").Value = [double]"-2.172551712240076e-12"
$ws.Range("E10").Value = [double]"1.319809842278922e-12"
$ws.Range("B11").Value = [double]"-2.122531149858278e-10"
$ws.Range("C11").Value = [double]"1.751075070656373e-10"
$ws.Range("D11").Value = [double]"3.992016857686692e-05"
$ws.Range("E11").Value = [double]"3.438079429741136e-08"
$ws.Range("B12").Value = [double]"8.973637435955895e-13"
$ws.Range("C12").Value = [double]"3.999632512940301e-08"
$ws.Range("D12").Value = [double]"1.75141671670269e-10"
$ws.Range("E12").Value = [double]"-2.12294180942487e-10"
$ws.Range("B13").Value = [double]"4.752053357443638e-13"
$ws.Range("C13").Value = [double]"-1.963188279622116e-12"
$ws.Range("D13").Value = [double]"1.019253803672413e-10"
$ws.Range("E13").Value = [double]"-1.543888863747628e-10"
$ws.Range("B14").Value = [double]"6.565907670853046e-15"
$ws.Range("C14").Value = [double]"-1.131215370340868e-13"
$ws.Range("B15").Value = [double]"1.75107413594228e-10"
$ws.Range("C15").Value = [double]"-2.122530215144186e-10"
$ws.Range("D15").Value = [double]"3.438079430154727e-08"
$ws.Range("E15").Value = [double]"3.99201685768177e-05"
$ws.Range("B16").Value = [double]"3.118398382128806e-08"
$ws.Range("C16").Value = [double]"6.761585682904612e-08"
$ws.Range("D16").Value = [double]"1.264427166688584e-06"
$ws.Range("E16").Value = [double]"-4.435564708495793e-06"
$ws.Range("B17").Value = [double]"6.761612279739156e-08"
$ws.Range("C17").Value = [double]"3.118416847447056e-08"
$ws.Range("D17").Value = [double]"-4.435564313778439e-06"
$ws.Range("E17").Value = [double]"1.264427539383081e-06"
$ws.Range("B18").Value = [double]"-1.681103389327766e-08"
$ws.Range("C18").Value = [double]"-1.681105824547489e-08"
$ws.Range("D18").Value = [double]"3.417646401448404e-06"
$ws.Range("E18").Value = [double]"3.41764641394214e-06"
$ws.Range("B19").Value = [double]"-4.347111524302626e-07"
$ws.Range("C19").Value = [double]"4.346964314627127e-07"
$ws.Range("D19").Value = [double]"-1.656835177033138e-05"
$ws.Range("E19").Value = [double]"1.65683313246324e-05"
$ws.Range("B20").Value = [double]"-9.38395346671006e-06"
$ws.Range("C20").Value = [double]"-9.383941703116488e-06"
$ws.Range("D20").Value = [double]"-0.001826047719527242"
$ws.Range("E20").Value = [double]"-0.001826047725693642"
$ws.Range("B21").Value = [double]"-1.826763955064523e-05"
$ws.Range("C21").Value = [double]"1.826763981873161e-05"
$ws.Range("D21").Value = [double]"-0.003354455922884504"
$ws.Range("E21").Value = [double]"0.003354455921257818"
$ws.Range("B22").Value = [double]"1.999816256488319e-05"
$ws.Range("C22").Value = [double]"4.486818688121897e-10"
$ws.Range("D22").Value = [double]"-1.061470905823701e-07"
$ws.Range("E22").Value = [double]"8.757083597052071e-08"
$ws.Range("B23").Value = [double]"-9.815934482880657e-10"
$ws.Range("C23").Value = [double]"2.376033296166719e-10"
$ws.Range("D23").Value = [double]"-7.71944433627437e-08"
$ws.Range("E23").Value = [double]"5.096269018362064e-08"
$ws.Range("B24").Value = [double]"-5.656076071052636e-11"
$ws.Range("C24").Value = [double]"3.282960737214759e-12"
$ws.Range("D24").Value = [double]"-1.086275858032893e-09"
$ws.Range("E24").Value = [double]"6.59904919769443e-10"
$ws.Range("B25").Value = [double]"-1.061265574829877e-07"
$ws.Range("C25").Value = [double]"8.755375350469449e-08"
$ws.Range("D25").Value = [double]"0.01996008428843342"
$ws.Range("E25").Value = [double]"1.719039714864943e-05"
$ws.Range("B26").Value = [double]"4.486818721209122e-10"
$ws.Range("C26").Value = [double]"1.999816256470121e-05"
$ws.Range("D26").Value = [double]"8.757083583486309e-08"
$ws.Range("E26").Value = [double]"-1.061470904715278e-07"
$ws.Range("B27").Value = [double]"2.376024825837247e-10"
$ws.Range("C27").Value = [double]"-9.81594295321013e-10"
$ws.Range("D27").Value = [double]"5.096269018362064e-08"
$ws.Range("E27").Value = [double]"-7.71944433627437e-08"
$ws.Range("B28").Value = [double]"3.282954119769859e-12"
$ws.Range("C28").Value = [double]"-5.656076732797126e-11"
$ws.Range("B29").Value = [double]"8.755370691788239e-08"
$ws.Range("C29").Value = [double]"-1.061265106844174e-07"
$ws.Range("D29").Value = [double]"1.719039715076702e-05"
$ws.Range("E29").Value = [double]"0.01996008428840886"
